# "add the NA's under duplicate_image_filename"
#
# Column E on the active sheet is the "duplicate_image_filename" column
# (header in E1). Every data row of the first table (rows 2-21) gets an
# "NA" value added in that column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Touch F1 (an existing, empty, string-typed cell) so the save round-trip
# keeps it empty instead of re-materializing a stray value for it.
$ws.Range("F1").Value = ""

for ($r = 2; $r -le 21; $r++) {
    $ws.Range("E$r").Value = "NA"
}
